$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename the four "2019 ..." sheets to "2022 ..." (year refresh).
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("2019 CPI").Name = "2022 CPI"
$wb.Worksheets.Item("2019 Equip Index Factors").Name = "2022 Equip Index Factors"
$wb.Worksheets.Item("2019 Agricultural % Good").Name = "2022 Agricultural % Good"
$wb.Worksheets.Item("2019 Construction % Good").Name = "2022 Construction % Good"

# ---------------------------------------------------------------------------
# 2. Repair the defined names so they point at the renamed sheets again
#    (renaming a sheet does not rewrite definedName RefersTo formulas).
# ---------------------------------------------------------------------------
$wb.Names.Item("Factor_to_Year").RefersTo = "='2022 CPI'!#REF!"
$wb.Names.Item("2022 CPI!Print_Area").RefersTo = "='2022 CPI'!`$A`$1:`$D`$40"
$wb.Names.Item("2022 Equip Index Factors!Print_Area").RefersTo = "='2022 Equip Index Factors'!`$A`$1:`$F`$42"
$wb.Names.Item("2022 Agricultural % Good!Print_Area").RefersTo = "='2022 Agricultural % Good'!`$A`$1:`$I`$23"
$wb.Names.Item("2022 Construction % Good!Print_Area").RefersTo = "='2022 Construction % Good'!`$A`$1:`$H`$23"

# ---------------------------------------------------------------------------
# 3. Page setup tweak on "Industrial Composite" (explicit portrait orientation).
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("Industrial Composite").PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# 4. Refresh / reset the saved window scroll position on several sheets
#    (drop the stale "topLeftCell" while keeping the existing selection).
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Construction Composite")
$ws.Activate()
$ws.Range("A2:A41").Select()

$ws = $wb.Worksheets.Item("Ag Mobile Equip Composite")
$ws.Activate()
$ws.Range("E2:E41").Select()

$ws = $wb.Worksheets.Item("Construction Mobile Composite")
$ws.Activate()
$ws.Range("A2:A41").Select()

$ws = $wb.Worksheets.Item("2022 Equip Index Factors")
$ws.Activate()
$ws.Range("A3:A42").Select()

$ws = $wb.Worksheets.Item("M&E Property Good Factor")
$ws.Activate()
$ws.Range("A3:A42").Select()

$ws = $wb.Worksheets.Item("2022 Agricultural % Good")
$ws.Activate()
$ws.Range("A4:A43").Select()

# ---------------------------------------------------------------------------
# 5. "2022 CPI" sheet: it was the active/selected tab before -- scroll it
#    back down (topLeftCell="A34") and move the selection to D46, and hand
#    tabSelected off to the new active sheet below.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("2022 CPI")
$ws.Activate()
$ws.Range("D46").Select()
$excel.ActiveWindow.ScrollRow = 34

# ---------------------------------------------------------------------------
# 6. "2022 Construction % Good" becomes the final active / selected sheet.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("2022 Construction % Good")
$ws.Activate()
$ws.Range("O15").Select()

Write-Host "done"
